# update: missing sessions from classes 1,2,3
#
# 1) The "Date Placeholder" field shown on every slide layout and on the
#    slide master is a fixed (non auto-updating) date that was typed as
#    18/10/2022 and needs to become 19/10/2022.
# 2) The rotated label next to the backlog chart on slide 9 needs its
#    wording tweaked from "reducción" to "manejo".

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "19/10/2022"

# Slide master's own date placeholder.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout has its own copy of the date placeholder too.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 9: rename the rotated "Indice de reduccion de backlog" textbox
# (nested inside a group) to "Indice de manejo de backlog".
$slide9 = $p.Slides.Item(9)
for ($i = 1; $i -le $slide9.Shapes.Count; $i++) {
    $topShape = $slide9.Shapes.Item($i)
    if ($topShape.Type -eq 6) {
        $items = $topShape.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $inner = $items.Item($j)
            if ($inner.Name -eq "TextBox 22") {
                $inner.TextFrame.TextRange.Text = "Índice de manejo de backlog"
            }
        }
    }
}
